# "Time recording log.xlsx" - log two more work sessions on the "Nädal 3"
# (Week 3) sheet: finish task 5 (row 11) and start task 6 (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 3")

# --- Row 11 (task 5, "Kodutöö 3"): record the stop time, duration,
#     activity and comment, and mark it complete ---
$ws.Range("D11").Value = 44005.962500000001
$ws.Range("F11").Value = 31
$ws.Range("G11").Value = "Kodutöö 3"
$ws.Range("H11").Value = "p. 8"
$ws.Range("J11").Value = "x"

# --- Row 12 (task 6, "Kodutöö 3"): record the start date/time and activity ---
$ws.Range("B12").Value = 43877
$ws.Range("C12").Value = 0.64027777777777783
$ws.Range("G12").Value = "Kodutöö 3"

# Columns G and H were manually narrowed (no longer auto-fit), and column D's
# auto-fit width now matches column C's (both show only a time value).
$ws.Columns.Item(3).ColumnWidth = 4.666666666666667
$ws.Columns.Item(4).ColumnWidth = 4.666666666666667
$ws.Columns.Item(7).ColumnWidth = 25.666666666666668
$ws.Columns.Item(8).ColumnWidth = 29.666666666666668

# Selection ends up on L7 after the edits
$ws.Range("L7").Select()
